# Commit: Change "Device" to "Apparatus" in excel form, simulink, function name
#
# 1. Rename worksheet "Device" -> "Apparatus"
# 2. Update the sheet's summary text, header row text/formatting
# 3. Make "Apparatus" the active/selected sheet (was "Advance")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# --- Rename the sheet ---
$ws.Name = "Apparatus"

# --- Update description text in A1 (device -> apparatus wording) ---
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# --- Update header row (row 3): new wording + bold formatting to match ---
# (A1 already carries the bold "section title" style; copy that formatting
#  onto the header cells after setting their text so the values aren't
#  clobbered by the paste.)
$ws.Range("A3").Value = "Bus number"
$ws.Range("B3").Value = "Type"
$ws.Range("C3").Value = "Parameters"

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Widen column A slightly to fit the new "Bus number" header ---
$ws.Columns.Item(1).ColumnWidth = 12.285714285714285

# --- Make "Apparatus" the active sheet/tab (previously "Advance" was active) ---
$ws.Activate()
$ws.Range("C4").Select()
